$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '330.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '6.88%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '40.25'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '8.07%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.597'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '9.30%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08132'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.70%'
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '8.680'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '4.92%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.930'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.69%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.977'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.17%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9494'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '3.00%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1283'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '20.08%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1980'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '4.10%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09190'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.60%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03566'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '7.59%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09584'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.02%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001318'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-4.59%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006137'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '3.92%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.366'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.78%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.547'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '3.31%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3513'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.85%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.376'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '16.05%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1333'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '1.38%'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.79%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04429'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '1.70%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '2.34%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004321'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '1.08%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-14.26%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003994'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '37.67%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02512'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '15.66%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05247'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '4.45%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007806'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.95%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1431'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '5.75%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.008325'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-2.31%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002140'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '6.49%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01044'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '28.86%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006606'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.77%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.08%'
$ws.Range('B48').Value = 'BOLO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002836'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-13.95%'
$ws.Range('B49').Value = 'CoinbaseStockToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002403'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '139.75%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002103'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.08%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.08%'
